$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Miles" rule values: <340 -> <340.0 and >0 -> >0.0
$ws.Range("F2").Value = "<340.0"
for ($r = 3; $r -le 13; $r++) {
    $ws.Cells.Item($r, 6).Value = ">0.0"
}

# Update the selected cell in the sheet view
$ws.Range("I10").Select()
